$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 41671856
$ws.Range("J64").Value2 = 5978.1
$ws.Range("L64").Value2 = 5978.1
$ws.Range("N64").Value2 = -6474.1
$ws.Range("H67").Value2 = 41671856
$ws.Range("J67").Value2 = 5978.1
$ws.Range("L67").Value2 = 5978.1
$ws.Range("N67").Value2 = -7694.1
$ws.Range("H88").Value2 = 2741.7827
$ws.Range("I88").Value2 = 3895.3333
$ws.Range("J88").Value2 = 2334.647
$ws.Range("K88").Value2 = 3895.3333
$ws.Range("L88").Value2 = 2334.647
$ws.Range("M88").Value2 = -3489.3333
$ws.Range("N88").Value2 = -3146.647
$ws.Range("H91").Value2 = 2741.7827
$ws.Range("I91").Value2 = 3895.3333
$ws.Range("J91").Value2 = 2334.647
$ws.Range("K91").Value2 = 3895.3333
$ws.Range("L91").Value2 = 2334.647
$ws.Range("M91").Value2 = -2491.3333
$ws.Range("N91").Value2 = -5142.647
$ws.Range("H92").Value2 = 2917.0527
$ws.Range("I92").Value2 = 3899.6
$ws.Range("K92").Value2 = 3899.6
$ws.Range("M92").Value2 = -2651.6
$ws.Range("H103").Value2 = 510.16666
$ws.Range("I103").Value2 = 556.5
$ws.Range("K103").Value2 = 1669.5
$ws.Range("M103").Value2 = -1083.5
$ws.Range("H116").Value2 = 9000
$ws.Range("J116").Value2 = 10000
$ws.Range("L116").Value2 = 10000
$ws.Range("N116").Value2 = -16884
$ws.Range("H137").Value2 = 5630.9585
$ws.Range("I137").Value2 = 1285.0488
$ws.Range("K137").Value2 = 3855.1464
$ws.Range("M137").Value2 = -1305.1464
$ws.Range("H138").Value2 = 3191.6099
$ws.Range("J138").Value2 = 2199.2942
$ws.Range("L138").Value2 = 6597.882599999999
$ws.Range("N138").Value2 = -16877.8826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1627.7333
$ws.Range("I2").Value2 = 1441.6
$ws.Range("K2").Value2 = 1441.6
$ws.Range("M2").Value2 = -1328.6
$ws.Range("H9").Value2 = 400001
$ws.Range("J9").Value2 = 400001
$ws.Range("L9").Value2 = 400001
$ws.Range("N9").Value2 = -400341
$ws.Range("H20").Value2 = 400001
$ws.Range("J20").Value2 = 400001
$ws.Range("L20").Value2 = 400001
$ws.Range("N20").Value2 = -400541
$ws.Range("H37").Value2 = 12400
$ws.Range("I37").Value2 = 12400
$ws.Range("K37").Value2 = 12400
$ws.Range("M37").Value2 = -12127
$ws.Range("H44").Value2 = 59500
$ws.Range("J44").Value2 = 59500
$ws.Range("L44").Value2 = 59500
$ws.Range("N44").Value2 = -60476
$ws.Range("H46").Value2 = 4859
$ws.Range("J46").Value2 = 4622.9
$ws.Range("L46").Value2 = 4622.9
$ws.Range("N46").Value2 = -5260.9
$ws.Range("H61").Value2 = 858507.4399999999
$ws.Range("I61").Value2 = 1198.9667
$ws.Range("K61").Value2 = 1198.9667
$ws.Range("M61").Value2 = -986.9666999999999
$ws.Range("H74").Value2 = 20332.355
$ws.Range("I74").Value2 = 1530.52
$ws.Range("K74").Value2 = 1530.52
$ws.Range("M74").Value2 = -656.52
$ws.Range("H77").Value2 = 20332.355
$ws.Range("I77").Value2 = 1530.52
$ws.Range("K77").Value2 = 7652.6
$ws.Range("M77").Value2 = -3284.6
$ws.Range("H102").Value2 = 4090
$ws.Range("J102").Value2 = 4066.6667
$ws.Range("L102").Value2 = 4066.6667
$ws.Range("N102").Value2 = -7310.6667
$ws.Range("H116").Value2 = 1627.7333
$ws.Range("I116").Value2 = 1441.6
$ws.Range("K116").Value2 = 1441.6
$ws.Range("M116").Value2 = 852.4000000000001
$ws.Range("H126").Value2 = 10099.2
$ws.Range("I126").Value2 = 10099.2
$ws.Range("K126").Value2 = 30297.6
$ws.Range("M126").Value2 = -27827.6
$ws.Range("H132").Value2 = 2764209.5
$ws.Range("I132").Value2 = 2346
$ws.Range("J132").Value2 = 11049800
$ws.Range("K132").Value2 = 7038
$ws.Range("L132").Value2 = 33149400
$ws.Range("M132").Value2 = -4508
$ws.Range("N132").Value2 = -33154460
$ws.Range("H136").Value2 = 858507.4399999999
$ws.Range("I136").Value2 = 1198.9667
$ws.Range("K136").Value2 = 3596.9001
$ws.Range("M136").Value2 = -1046.9001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1627.7333
$ws.Range("I3").Value2 = 1441.6
$ws.Range("K3").Value2 = 1441.6
$ws.Range("M3").Value2 = -1327.6
$ws.Range("H94").Value2 = 2472.8333
$ws.Range("I94").Value2 = 1113.45
$ws.Range("J94").Value2 = 9269.75
$ws.Range("K94").Value2 = 1113.45
$ws.Range("L94").Value2 = 9269.75
$ws.Range("M94").Value2 = -662.45
$ws.Range("N94").Value2 = -10171.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value2 = 3833
$ws.Range("I25").Value2 = 3833
$ws.Range("K25").Value2 = 3833
$ws.Range("M25").Value2 = -3659
$ws.Range("H99").Value2 = 6877022.5
$ws.Range("I99").Value2 = 4454031
$ws.Range("K99").Value2 = 4454031
$ws.Range("M99").Value2 = -4452533
$ws.Range("H126").Value2 = 6877022.5
$ws.Range("I126").Value2 = 4454031
$ws.Range("K126").Value2 = 13362093
$ws.Range("M126").Value2 = -13359623
$ws.Range("H132").Value2 = 73533224
$ws.Range("I132").Value2 = 3702.889
$ws.Range("J132").Value2 = 205886350
$ws.Range("K132").Value2 = 11108.667
$ws.Range("L132").Value2 = 617659050
$ws.Range("M132").Value2 = -8578.667000000001
$ws.Range("N132").Value2 = -617664110
$ws.Range("H134").Value2 = 34489396
$ws.Range("I134").Value2 = 1471.0555
$ws.Range("K134").Value2 = 4413.166499999999
$ws.Range("M134").Value2 = -1878.166499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 19.647058
$ws.Range("I12").Value2 = 12
$ws.Range("K12").Value2 = 36
$ws.Range("M12").Value2 = 137
$ws.Range("H87").Value2 = 1252.8
$ws.Range("I87").Value2 = 1252.8
$ws.Range("K87").Value2 = 3758.4
$ws.Range("M87").Value2 = -2510.4
$ws.Range("H90").Value2 = 1252.8
$ws.Range("I90").Value2 = 1252.8
$ws.Range("K90").Value2 = 11275.2
$ws.Range("M90").Value2 = -5035.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 6524.35
$ws.Range("I70").Value2 = 5089.5835
$ws.Range("K70").Value2 = 5089.5835
$ws.Range("M70").Value2 = -4819.5835
$ws.Range("H73").Value2 = 6524.35
$ws.Range("I73").Value2 = 5089.5835
$ws.Range("K73").Value2 = 5089.5835
$ws.Range("M73").Value2 = -4153.5835
$ws.Range("H113").Value2 = 2855.5
$ws.Range("J113").Value2 = 2846.8
$ws.Range("L113").Value2 = 2846.8
$ws.Range("N113").Value2 = -7186.8
$ws.Range("H122").Value2 = 1261409.6
$ws.Range("I122").Value2 = 1547573.1
$ws.Range("J122").Value2 = 2290
$ws.Range("K122").Value2 = 4642719.300000001
$ws.Range("L122").Value2 = 6870
$ws.Range("M122").Value2 = -4640269.300000001
$ws.Range("N122").Value2 = -11770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2807917.5
$ws.Range("J7").Value2 = 6666.3335
$ws.Range("L7").Value2 = 6666.3335
$ws.Range("N7").Value2 = -6890.3335
$ws.Range("H22").Value2 = 100001430
$ws.Range("I22").Value2 = 1499.5
$ws.Range("J22").Value2 = 125001416
$ws.Range("K22").Value2 = 1499.5
$ws.Range("L22").Value2 = 125001416
$ws.Range("M22").Value2 = -1204.5
$ws.Range("N22").Value2 = -125002006
$ws.Range("H27").Value2 = 100001430
$ws.Range("I27").Value2 = 1499.5
$ws.Range("J27").Value2 = 125001416
$ws.Range("K27").Value2 = 1499.5
$ws.Range("L27").Value2 = 125001416
$ws.Range("M27").Value2 = -1392.5
$ws.Range("N27").Value2 = -125001630
$ws.Range("H93").Value2 = 142873390
$ws.Range("I93").Value2 = 333350600
$ws.Range("J93").Value2 = 15498.5
$ws.Range("K93").Value2 = 333350600
$ws.Range("L93").Value2 = 15498.5
$ws.Range("M93").Value2 = -333349352
$ws.Range("N93").Value2 = -17994.5
$ws.Range("H126").Value2 = 2807917.5
$ws.Range("J126").Value2 = 6666.3335
$ws.Range("L126").Value2 = 19999.0005
$ws.Range("N126").Value2 = -24939.0005
$ws.Range("H132").Value2 = 1592758.9
$ws.Range("I132").Value2 = 3841.7932
$ws.Range("J132").Value2 = 4664665
$ws.Range("K132").Value2 = 11525.3796
$ws.Range("L132").Value2 = 13993995
$ws.Range("M132").Value2 = -8995.3796
$ws.Range("N132").Value2 = -13999055
$ws.Range("H136").Value2 = 1076480.5
$ws.Range("I136").Value2 = 15542.571
$ws.Range("K136").Value2 = 46627.713
$ws.Range("M136").Value2 = -44077.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value2 = 5641000
$ws.Range("I3").Value2 = 9366666
$ws.Range("J3").Value2 = 52500
$ws.Range("K3").Value2 = 9366666
$ws.Range("L3").Value2 = 52500
$ws.Range("M3").Value2 = -9366552
$ws.Range("N3").Value2 = -52728
$ws.Range("H96").Value2 = 1811.75
$ws.Range("I96").Value2 = 1598.75
$ws.Range("J96").Value2 = 1918.25
$ws.Range("K96").Value2 = 1598.75
$ws.Range("L96").Value2 = 1918.25
$ws.Range("M96").Value2 = -225.75
$ws.Range("N96").Value2 = -4664.25
$ws.Range("H107").Value2 = 1065.3077
$ws.Range("I107").Value2 = 1065.3077
$ws.Range("K107").Value2 = 3195.9231
$ws.Range("M107").Value2 = -1275.9231
$ws.Range("H122").Value2 = 422361.78
$ws.Range("I122").Value2 = 586859.5
$ws.Range("K122").Value2 = 1760578.5
$ws.Range("M122").Value2 = -1758128.5
